$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the two new time-log entries (rows 30 and 31) ---
# Copy formatting (borders / date number format) down from the row above
# so the new rows look exactly like the existing entries (row 28/29 block).
$ws.Range("A28:F29").Copy()
$ws.Range("A30:F31").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Row 30 - Viki, 2024-02-27 (serial 45349), 8:00 -> 9:35
$ws.Range("A30").Value = "Viki"
$ws.Range("B30").Value = 45349
$ws.Range("C30").Formula = "=8"

# Row 31 - Aris, 2024-02-27 (serial 45349), 8:00 -> 9:35
$ws.Range("A31").Value = "Aris"
$ws.Range("B31").Value = 45349
$ws.Range("C31").Formula = "=8"

# D24:D31 all share the same "end time" formula; (re)writing the whole
# contiguous block at once turns it into one shared formula group
# (si="1"), exactly covering the now-extended range through D31.
$ws.Range("D24:D31").Formula = "=9+35/60"

# --- Update the view state to match where the author ended up ---
$ws.Range("D31").Select()
